$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay as literal text (matches source formatting)
$textCells = @("D5", "D6", "D8", "D16", "D17", "D25", "D26", "D36", "D37", "D41", "D45", "D48", "D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

# Apply cell value updates
$ws.Range("D2").Value = '27.560.59'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '1.660.63'
$ws.Range("E3").Value = '  -3.61%  '
$ws.Range("E4").Value = '  +0.67%  '
$ws.Range("D5").Value = '214.60'
$ws.Range("E5").Value = '  -1.35%  '
$ws.Range("D6").Value = '0.513'
$ws.Range("E6").Value = '  -1.71%  '
$ws.Range("E7").Value = '  +0.66%  '
$ws.Range("D8").Value = '23.31'
$ws.Range("E8").Value = '  -3.12%  '
$ws.Range("E9").Value = '  -2.66%  '
$ws.Range("E10").Value = '  -1.90%  '
$ws.Range("E11").Value = '  -2.40%  '
$ws.Range("D12").Value = '1.894.75'
$ws.Range("E12").Value = '  -3.62%  '
$ws.Range("D13").Value = '1.659.25'
$ws.Range("E13").Value = '  -3.66%  '
$ws.Range("E14").Value = '  -2.91%  '
$ws.Range("E15").Value = '  -3.11%  '
$ws.Range("D16").Value = '65.75'
$ws.Range("E16").Value = '  -3.31%  '
$ws.Range("D17").Value = '246.39'
$ws.Range("E17").Value = '  +1.61%  '
$ws.Range("D18").Value = '27.543.08'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").Value = '0.0₃0731'
$ws.Range("E19").Value = '  -2.46%  '
$ws.Range("E20").Value = '  -6.83%  '
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("E22").Value = '  -3.44%  '
$ws.Range("E23").Value = '  -4.38%  '
$ws.Range("E24").Value = '  -4.66%  '
$ws.Range("D25").Value = '145.96'
$ws.Range("E25").Value = '  -1.74%  '
$ws.Range("D26").Value = '7.16'
$ws.Range("E26").Value = '  -5.39%  '
$ws.Range("E27").Value = '  -2.47%  '
$ws.Range("E29").Value = '  -2.15%  '
$ws.Range("E30").Value = '  +4.64%  '
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("E32").Value = '  -3.13%  '
$ws.Range("D33").Value = '1.441.89'
$ws.Range("E33").Value = '  -7.05%  '
$ws.Range("E34").Value = '  -5.40%  '
$ws.Range("E35").Value = '  -8.35%  '
$ws.Range("D36").Value = '2.37'
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("D37").Value = '0.929'
$ws.Range("E37").Value = '  -4.20%  '
$ws.Range("E38").Value = '  -6.06%  '
$ws.Range("E39").Value = '  -2.90%  '
$ws.Range("E40").Value = '  -2.02%  '
$ws.Range("D41").Value = '69.03'
$ws.Range("E41").Value = '  -3.53%  '
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("E43").Value = '  -7.67%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").Value = '2.20'
$ws.Range("E45").Value = '  -3.51%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.803.22'
$ws.Range("E46").Value = '  -3.00%  '
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("D48").Value = '88.65'
$ws.Range("E48").Value = '  -3.80%  '
$ws.Range("E49").Value = '  -2.80%  '
$ws.Range("E50").Value = '  -4.54%  '
$ws.Range("D51").Value = '7.81'
$ws.Range("E51").Value = '  -6.69%  '
